# Generate Report for handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the first data row on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-17 10:14:33"
$wsZhCn.Range("G2").Value = "2016-01-17 10:15:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-17 10:14:44"
$wsDeDe.Range("G2").Value = "2016-01-17 10:16:03"
